$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("A2").Value = "30 Oct 2025, 11:53 AM"
$ws = $wb.Worksheets.Item("1 Month Performance")
$ws.Range("C3").Value = 81.8485
$ws.Range("C4").Value = 78.11839999999999
$ws.Range("C5").Value = 66.8811
$ws.Range("C6").Value = 60.9253
$ws.Range("C7").Value = 59.5828
$ws.Range("C9").Value = 54.0314
$ws.Range("C10").Value = 44.756
$ws.Range("B11").Value = "MTARTECH"
$ws.Range("C11").Value = 40.6503
$ws.Range("B12").Value = "V2RETAIL"
$ws.Range("C12").Value = 40.6422
$ws.Range("C14").Value = 38.9017
$ws.Range("B15").Value = "TVSSRICHAK"
$ws.Range("C15").Value = 38.4595
$ws.Range("B16").Value = "SHAREINDIA"
$ws.Range("C16").Value = 37.6324
$ws.Range("C17").Value = 36.5272
$ws.Range("C18").Value = 35.6671
$ws.Range("B19").Value = "MEGASOFT"
$ws.Range("C19").Value = 34.8668
$ws.Range("B20").Value = "SAMMAANCAP"
$ws.Range("C20").Value = 34.538
$ws.Range("C22").Value = 33.6829
$ws.Range("C23").Value = 33.1681
$ws.Range("C24").Value = 32.0291
$ws.Range("C25").Value = 31.7901
$ws.Range("B27").Value = "ORIENTTECH"
$ws.Range("C27").Value = 29.0468
$ws.Range("B28").Value = "TARACHAND"
$ws.Range("C28").Value = 28.7348
$ws.Range("C31").Value = 27.4089
$ws.Range("C33").Value = 26.6572
$ws.Range("C35").Value = 25.8332
$ws.Range("C36").Value = 25.4825
$ws.Range("C37").Value = 24.7282
$ws.Range("C38").Value = 24.4943
$ws.Range("C39").Value = 24.3214
$ws.Range("C40").Value = 24.0285
$ws.Range("C45").Value = 23.3896
$ws.Range("B47").Value = "INDIANB"
$ws.Range("C47").Value = 22.8395
$ws.Range("B48").Value = "MARINE"
$ws.Range("C48").Value = 22.8001
$ws.Range("C49").Value = 22.7162
$ws.Range("C50").Value = 22.3719
$ws.Range("B51").Value = "KERNEX"
$ws.Range("C51").Value = 22.2922
$ws.Range("B52").Value = "IIFL"
$ws.Range("C52").Value = 22.2745
$ws.Range("C53").Value = 22.2003
$ws.Range("C54").Value = 22.0326
$ws.Range("C55").Value = 21.8222
$ws.Range("B56").Value = "GUJTHEM"
$ws.Range("C56").Value = 21.3115
$ws.Range("B57").Value = "PRIVISCL"
$ws.Range("C57").Value = 21.2558
$ws.Range("C60").Value = 20.503
$ws.Range("B62").Value = "INDRAMEDCO"
$ws.Range("C62").Value = 20.3112
$ws.Range("B63").Value = "SHRIRAMFIN"
$ws.Range("C63").Value = 20.2071
$ws.Range("C65").Value = 19.8247
$ws.Range("C67").Value = 19.5198
$ws.Range("C68").Value = 19.4463
$ws.Range("C69").Value = 19.1829
$ws.Range("C70").Value = 19.0664
$ws.Range("C74").Value = 18.4784
$ws.Range("B76").Value = "ACUTAAS"
$ws.Range("C76").Value = 18.2516
$ws = $wb.Worksheets.Item("distance from Dma50")
$ws.Range("C2").Value = 9.8567
$ws.Range("C3").Value = 7.4596
$ws.Range("C4").Value = 6.6904
$ws.Range("C5").Value = 5.4401
$ws.Range("C6").Value = 5.4036
$ws.Range("C7").Value = 5.1949
$ws.Range("C8").Value = 4.5158
$ws.Range("C9").Value = 4.4343
$ws.Range("C10").Value = 3.9279
$ws.Range("C11").Value = 3.8634
$ws.Range("C12").Value = 3.5077
$ws.Range("C13").Value = 3.4578
$ws.Range("C14").Value = 3.1617
$ws.Range("C15").Value = 3.1306
$ws.Range("C16").Value = 3.0522
$ws.Range("C17").Value = 2.906
$ws.Range("C18").Value = 2.9027
$ws.Range("C19").Value = 2.817
$ws.Range("C20").Value = 2.4724
$ws.Range("C21").Value = 2.3847
$ws.Range("C22").Value = 1.4142
$ws.Range("C23").Value = 1.4117
$ws.Range("C24").Value = 1.3206
$ws.Range("C25").Value = 1.0447
$ws.Range("C26").Value = 1.0017
$ws.Range("C27").Value = 0.8844
$ws.Range("C28").Value = 0.6002999999999999
$ws.Range("C29").Value = 0.3955
$ws.Range("C30").Value = -2.0867
